$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New entry for row 25 (08 June 2023 update)
$ws.Range("A25").Value = Get-Date -Year 2023 -Month 6 -Day 6 -Hour 0 -Minute 0 -Second 0
$ws.Range("B25").NumberFormat = $ws.Range("B24").NumberFormat
$ws.Range("B25").Value = 39000
$ws.Range("C25").Value = 0
$ws.Range("E25").Value = "uang persembahan - reguler"
$ws.Range("F25").Value = "yofandi"

# Update view: scroll/selection position (topLeftCell A7, active cell B27)
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B27").Select()
